$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view size (workbook.xml bookViews/workbookView) ---
$excel.Width = 28800
$excel.Height = 12000

# --- Fix up two stray Department lookups (CANCELLED / INVOICED rows) ---
$ws.Range("B2").Value = "CANCELLED"
$ws.Range("B7").Value = "INVOICED"

# --- Rewrite the tail of the OrderStatus/Department table (rows 269-277) ---
# "FRAME REVEIVED" is dropped from the list, shifting everything after it up
# by one logical slot, and the department assignments are refreshed.
$ws.Range("A269").Value = "SEND TO LENSWARE"
$ws.Range("B269").Value = "CS"

$ws.Range("A270").Value = "BUCODE CHANGED"
$ws.Range("B270").Value = "FITT"

$ws.Range("A271").Value = "ISSUED TO FITTING (A2)-GRG"
$ws.Range("B271").Value = "FITT"

$ws.Range("A272").Value = "ISSUED TO FINAL FITTING-QC (A2)- GRG"
$ws.Range("B272").Value = "FITT"

$ws.Range("A273").Value = "ISSUED TO FITTING-QC-IN (A14)-GRG"
$ws.Range("B273").Value = "FITT"

$ws.Range("A274").Value = "ISSUED TO FITTING (A14)-GRG"
$ws.Range("B274").Value = "FITT"

$ws.Range("A275").Value = "CUSTOMER COPY PRINTED."
$ws.Range("B275").Value = "CS"

$ws.Range("A276").Value = "PRODUCTION COPY PRINTED."
$ws.Range("B276").Value = "CS"

$ws.Range("A277").Value = "ISSUED TO MOUNTIING (A2)-GRG"
$ws.Range("B277").Value = "FITT"
$ws.Range("C277").Value = "Target"

# --- New departments / logic for movement to TC (rows 278-289) ---
$ws.Range("A278").Value = "ISSUED TO MOUNT REPROCESS(A2)-GRG"
$ws.Range("B278").Value = "FITT"

$ws.Range("A279").Value = "ISSUED TO MOUNT REPROCESS(A14)-GRG"
$ws.Range("B279").Value = "FITT"

$ws.Range("A280").Value = "ISSUED TO MOUNTING (A2)-GRG"
$ws.Range("B280").Value = "FITT"

$ws.Range("A281").Value = "ISSUED TO MOUNTING (A14)-GRG"
$ws.Range("B281").Value = "FITT"

$ws.Range("A282").Value = "ISSUED TO SURFACING-DS REPROCESS(A14-SF)-GRG"
$ws.Range("B282").Value = "DS"

$ws.Range("A283").Value = "ISSUED TO TINT REPROCESS (A2)-GRG"
$ws.Range("B283").Value = "TINT"

$ws.Range("A284").Value = "RETURNED TO MOUNT REPROCESS(A2)-GRG"
$ws.Range("B284").Value = "FITT"

$ws.Range("A285").Value = "RETURNED TO MOUNT REPROCESS(A14)-GRG"
$ws.Range("B285").Value = "FITT"

$ws.Range("A286").Value = "RETURNED TO MOUNTING (A2)-GRG"
$ws.Range("B286").Value = "FITT"

$ws.Range("A287").Value = "RETURNED TO MOUNTING (A14)-GRG"
$ws.Range("B287").Value = "FITT"

$ws.Range("A288").Value = "RETURNED TO SURFACING-DS REPROCESS(A14-SF)-GRG"
$ws.Range("B288").Value = "DS"

$ws.Range("A289").Value = "RETURNED TO TINT REPROCESS (A2)-GRG"
$ws.Range("B289").Value = "TINT"

# --- Column A width grows to fit the longer new department names ---
# (best-fit width target ~50.43 chars; this runtime quantizes ColumnWidth,
# so use the nearest value that round-trips closest to that target)
$ws.Columns.Item(1).ColumnWidth = 49.6

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / scroll position ---
$ws.Range("B7").Select()
